$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:R47").AutoFilter()

$sortRange = $ws.Range("A1:R47")
$key1 = $ws.Range("G1:G47")
$sortRange.Sort($key1, 2, $null, $null, 1, $null, 1, 1)
